# Add 2022-Q4 data
#
# Starting layout:  总计 | 2022-Q3
# Target layout:     总计 | 2022-Q4 | 2022-Q3
#
# The sheet that currently holds the Q3 fund-holdings table is repurposed to hold the
# new Q4 numbers (it keeps its tab position right after 总计), and a fresh sheet named
# "2022-Q3" is appended after it, preserving the old Q3 numbers that used to live there.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets("总计")
$q3 = $wb.Worksheets("2022-Q3")

# --- 1) Spin up the new archive tab (placed right after the existing one) and copy the
#        existing Q3 table + formatting into it verbatim, before we overwrite the original. ---
$q3Archive = $wb.Worksheets.Add($null, $q3)

$q3.Range("B1:H1").Copy($q3Archive.Range("B1"))
$q3.Range("A2:H2").Copy($q3Archive.Range("A2"))

# --- 2) Repurpose the original Q3 sheet into the Q4 sheet: rename the tab and replace
#        its data with the new quarter's numbers. Rename the original sheet away from
#        "2022-Q3" before claiming that name for the archive copy, to avoid a name clash. ---
$q3.Name = "2022-Q4"
$q4 = $q3
$q3Archive.Name = "2022-Q3"

# Match the header/A2 styling used on the 总计 sheet (bold + border look already present there).
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)

$q4.Range("A2").Value = 0
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("B2").Value = "690003"
$q4.Range("C2").Value = "民生加银精选混合"
$q4.Range("D2").Value = "0.50"
$q4.Range("E2").Value = "88.58"
$q4.Range("F2").Value = "6.53"
$q4.Range("G2").Value = "0.0326"
$q4.Range("H2").Value = 2

# --- 3) Update the 总计 sheet: relabel row 2 as 2022-Q4 and append a restored 2022-Q3 row. ---
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.03

$summary.Range("B2").Value = "2022-Q4"

Write-Output "2022-Q4 sheet added"
